# Weekly refresh of the "Pera asiática" price rows: the D (Fecha), L (Calidad),
# M (Volumen), N/O/P (Precio mínimo/máximo/promedio ponderado), Q (Unidad de
# comercialización), R (Origen), S (Precio $/Kg) and T (Kg / unidad) columns
# are re-populated with the latest weekly data for rows 2-6 and 8-11
# (row 7 is unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{
    2  = @{ D = 44601; L = 'Primera';     M = 30;  N = 28000; O = 28000; P = 28000; Q = '$/caja 18 kilos granel';   R = "Región de O'Higgins"; S = 1556; T = 18 }
    3  = @{ D = 44495; L = 'Primera';     M = 50;  N = 24000; O = 24000; P = 24000; Q = '$/bandeja 10 kilos';        R = 'China';                S = 2400; T = 10 }
    4  = @{ D = 44427; L = 'Primera';     M = 55;  N = 7000;  O = 7000;  P = 7000;  Q = '$/caja 15 kilos granel';   R = "Región de O'Higgins"; S = 467;  T = 15 }
    5  = @{ D = 44208; L = 'Especial';    M = 70;  N = 24000; O = 24000; P = 24000; Q = '$/caja 15 kilos granel';   R = "Región de O'Higgins"; S = 1600; T = 15 }
    6  = @{ D = 44392; L = 'Especial';    M = 500; N = 7000;  O = 7000;  P = 7000;  Q = '$/bandeja 8 kilos';         R = "Región de O'Higgins"; S = 875;  T = 8  }
    8  = @{ D = 44418; L = 'Especial';    M = 100; N = 8000;  O = 8000;  P = 8000;  Q = '$/caja 15 kilos granel';   R = "Región de O'Higgins"; S = 533;  T = 15 }
    9  = @{ D = 44511; L = 'Primera';     M = 15;  N = 22000; O = 22000; P = 22000; Q = '$/caja 15 kilos granel';   R = "Región de O'Higgins"; S = 1467; T = 15 }
    10 = @{ D = 44264; L = 'Calibre 100'; M = 50;  N = 20000; O = 20000; P = 20000; Q = '$/caja 18 kilos embalada'; R = "Región de O'Higgins"; S = 1111; T = 18 }
    11 = @{ D = 44217; L = 'Primera';     M = 55;  N = 18000; O = 18000; P = 18000; Q = '$/caja 18 kilos granel';   R = "Región de O'Higgins"; S = 1000; T = 18 }
}

foreach ($r in $rows.Keys) {
    $row = $rows[$r]
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = $row.Q
    $ws.Cells.Item($r, 18).Value = $row.R
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = $row.T
}
